$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 ("OBI_DEF"), matching the style of the existing header row (E1)
$ws.Range("F1").Value = "OBI_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new OBI_DEF column (F) with the OBI definition values for each row
$ws.Range("F2").Value = "['A material entity is an independent continuant that has some portion of matter as proper or improper continuant part. [BFO]']"
$ws.Range("F3").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]"
$ws.Range("F4").Value = "['p is a process boundary if p is a temporal part of a process and p has no proper temporal parts. [BFO]']"
$ws.Range("F5").Value = "['A function is a disposition that exists in virtue of the bearer’s physical make-up and this physical make-up is something the bearer possesses because it came into being, either through evolution (in the case of natural biological entities) or through intentional design (in the case of artifacts), in order to realize processes of a certain sort. [BFO]']"
$ws.Range("F6").Value = "['A quality is a specifically dependent continuant that, in contrast to roles and dispositions, does not require any further process in order to be realized. [BFO]']"
$ws.Range("F7").Value = "['B is a process profile if there is some process c such that b is process profile of c. B is process profile of c holds when b is a proper occurrent part of c and there is some proper occurrent part d of c which has no parts in common with b and is mutually dependent on b and is such that b, c and d occupy the same temporal region. [BFO]']"
$ws.Range("F8").Value = "['An occurrent is an entity that unfolds itself in time or it is the instantaneous boundary of such an entity (for example a beginning or an ending) or it is a temporal or spatiotemporal region which such an entity occupies a temporal region or occupies a spatiotemporal region. [BFO]']"
$ws.Range("F9").Value = "['To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]', 'To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])']"
$ws.Range("F10").Value = "['A continuant is an entity that persists, endures, or continues to exist through time while maintaining its identity. [BFO]']"
$ws.Range("F11").Value = "['b is an independent continuant if b is a continuant which is such that there is no c and no t such that b s-depends_on c at t. [BFO]']"
$ws.Range("F12").Value = "['A spatial region is a continuant entity that is a continuant part of space R as defined relative to some frame R. [BFO]']"
$ws.Range("F13").Value = "['An entity is anything that exists or has existed or will exist. [BFO]']"
$ws.Range("F14").Value = "['A temporal region is an occurrent entity that is part of time as defined relative to some reference frame. [BFO]']"
$ws.Range("F15").Value = "['A spatiotemporal region is an occurrent entity that is part of spacetime. [BFO]']"
$ws.Range("F16").Value = "['b is a specifically dependent continuant if b is a continuant and there is some independent continuant c which is not a spatial region and which is such that b specifically depends on c at every time t during the course of b’s existence. [BFO]']"
$ws.Range("F17").Value = "['b is a generically dependent continuant if b is a continuant that generically depends on one or more other entities. [BFO]']"
$ws.Range("F18").Value = "['B is a role means: b is a realizable entity and b exists because there is some single bearer that is in some special physical, social, or institutional set of circumstances in which this bearer does not have to be and b is not such that, if it ceases to exist, then the physical make-up of the bearer is thereby changed. [BFO]']"
$ws.Range("F19").Value = "['B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]']"
$ws.Range("F20").Value = "['B is a fiat object part if b is a material entity which is such that for all times t, if b exists at t then there is some object c such that b is a proper continuant part of c at t and c is demarcated from the remainder of c by a two-dimensional continuant fiat boundary. [BFO]']"
$ws.Range("F21").Value = "['B is a site means: b is a three-dimensional immaterial entity that is (partially or wholly) bounded by a material entity or it is a three-dimensional immaterial part thereof. [BFO]']"
$ws.Range("F22").Value = "['B is an object means: b is a material entity which manifests causal unity of one or other of the types causal unities and is of a type (a material universal) instances of which are maximal relative to this criterion of causal unity.\n\nTo say that b is causally unified means: b is a material entity which is such that its material parts are tied together in such a way that, in environments typical for entities of the type in question,if c, a continuant part of b that is in the interior of b at t, is larger than a certain threshold size (which will be determined differently from case to case, depending on factors such as porosity of external cover) and is moved in space to be at t at a location on the exterior of the spatial region that had been occupied by b at t, then either b’s other parts will be moved in coordinated fashion or b will be damaged (be affected, for example, by breakage or tearing) in the interval between t and t. [BFO]']"
$ws.Range("F23").Value = "['B is an object aggregate means: b is a material entity consisting exactly of a plurality of objects as member parts at all times at which b exists. [BFO]']"
$ws.Range("F24").Value = "['B is a relational quality if for some independent continuants c, d and for some time t: b is quality of c at t and b is quality of d at t. [BFO]']"
$ws.Range("F25").Value = "['A history is a process that is the sum of the totality of processes taking place in the spatiotemporal region occupied by a material entity or site, including processes on the surface of the entity or within the cavities to which it serves as host. [BFO]']"
$ws.Range("F26").Value = "['A zero-dimensional temporal region is a temporal region that is without extent. [BFO]']"
$ws.Range("F27").Value = "['A one-dimensional temporal region is a temporal region that is extended. [BFO]']"
$ws.Range("F28").Value = "[]"
$ws.Range("F29").Value = "['B is a continuant fiat boundary if b is an immaterial entity that is of zero, one or two dimensions and does not include a spatial region as part. [BFO]']"
$ws.Range("F30").Value = "['A two-dimensional spatial region is a spatial region that is of two dimensions. [BFO]']"
$ws.Range("F31").Value = "['A three-dimensional spatial region is a spatial region that is of three dimensions. [BFO]']"
$ws.Range("F32").Value = "['A zero-dimensional spatial region is a point in space. [BFO]']"
$ws.Range("F33").Value = "['A one-dimensional spatial region is a line or aggregate of lines stretching from one point in space to another. [BFO]']"
$ws.Range("F34").Value = "['A one-dimensional continuant fiat boundary is a continuous fiat line whose location is defined in relation to some material entity.  [BFO]']"
$ws.Range("F35").Value = "['A two-dimensional continuant fiat boundary (surface) is a self-connected fiat surface whose location is defined in relation to some material entity. [BFO]']"
$ws.Range("F36").Value = "['A zero-dimensional continuant fiat boundary is a fiat point whose location is defined in relation to some material entity. [BFO]']"
